$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "Custos_Iniciativas" -> "Custos"
# ---------------------------------------------------------------------------
$wsCustos = $wb.Worksheets.Item("Custos_Iniciativas")
$wsCustos.Name = "Custos"

# ---------------------------------------------------------------------------
# 2) Parametros sheet: update absenteeism parameter values
# ---------------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Parametros")

$wsParam.Range("C4").Value = 0.1
$wsParam.Range("D4").Value = 0.001

$wsParam.Range("C5").Value = 0.1
$wsParam.Range("D5").Value = 0.001
$wsParam.Range("C5").Style = "Normal"
$wsParam.Range("D5").Style = "Normal"

$wsParam.Range("C6").Value = 10

$wsParam.Range("C7").Value = 20

$wsParam.Range("C8").Value = 0.05
$wsParam.Range("D8").Value = 0.001

$wsParam.Range("C9").Value = 0.05
$wsParam.Range("D9").Value = 0.001
$wsParam.Range("C9").Style = "Normal"
$wsParam.Range("D9").Style = "Normal"

$wsParam.Range("C12").Value = 0.1
$wsParam.Range("D12").Value = 0.001

$wsParam.Range("C13").Value = 0.05
$wsParam.Range("D13").Value = 0.001
$wsParam.Range("C13").Style = "Normal"
$wsParam.Range("D13").Style = "Normal"

$wsParam.Range("C16").Value = 0.1
$wsParam.Range("D16").Value = 0.001

$wsParam.Range("C17").Value = 0.05
$wsParam.Range("D17").Value = 0.001

# ---------------------------------------------------------------------------
# 3) Custos sheet: rename headers, update totals, add new rows for
#    "SemIniciativa" and "TodasIniciativas" scenarios
# ---------------------------------------------------------------------------
$wsCustos.Range("A1").Value = "Cenario"
$wsCustos.Range("D1").Value = "CustoTotal"

$wsCustos.Range("D2").Value = 50000
$wsCustos.Range("D3").Value = 50000
$wsCustos.Range("D4").Value = 50000
$wsCustos.Range("D5").Value = 50000
$wsCustos.Range("D6").Value = 50000

$wsCustos.Range("A12").Value = "SemIniciativa"
$wsCustos.Range("B12").Value = "Custo Total"
$wsCustos.Range("C12").Value = 2017
$wsCustos.Range("D12").Value = 50

$wsCustos.Range("A13").Value = "SemIniciativa"
$wsCustos.Range("B13").Value = "Custo Total"
$wsCustos.Range("C13").Value = 2018
$wsCustos.Range("D13").Value = 20

$wsCustos.Range("A14").Value = "SemIniciativa"
$wsCustos.Range("B14").Value = "Custo Total"
$wsCustos.Range("C14").Value = 2019
$wsCustos.Range("D14").Value = 20

$wsCustos.Range("A15").Value = "SemIniciativa"
$wsCustos.Range("B15").Value = "Custo Total"
$wsCustos.Range("C15").Value = 2020
$wsCustos.Range("D15").Value = 20

$wsCustos.Range("A16").Value = "SemIniciativa"
$wsCustos.Range("B16").Value = "Custo Total"
$wsCustos.Range("C16").Value = 2021
$wsCustos.Range("D16").Value = 20

$wsCustos.Range("A17").Value = "TodasIniciativas"
$wsCustos.Range("A17").Font.Bold = $true
$wsCustos.Range("B17").Value = "Custo Total"
$wsCustos.Range("C17").Value = 2017
$wsCustos.Range("D17").Value = 50

$wsCustos.Range("A18").Value = "TodasIniciativas"
$wsCustos.Range("A18").Font.Bold = $true
$wsCustos.Range("B18").Value = "Custo Total"
$wsCustos.Range("C18").Value = 2018
$wsCustos.Range("D18").Value = 20

$wsCustos.Range("A19").Value = "TodasIniciativas"
$wsCustos.Range("A19").Font.Bold = $true
$wsCustos.Range("B19").Value = "Custo Total"
$wsCustos.Range("C19").Value = 2019
$wsCustos.Range("D19").Value = 20

$wsCustos.Range("A20").Value = "TodasIniciativas"
$wsCustos.Range("A20").Font.Bold = $true
$wsCustos.Range("B20").Value = "Custo Total"
$wsCustos.Range("C20").Value = 2020
$wsCustos.Range("D20").Value = 20

$wsCustos.Range("A21").Value = "TodasIniciativas"
$wsCustos.Range("A21").Font.Bold = $true
$wsCustos.Range("B21").Value = "Custo Total"
$wsCustos.Range("C21").Value = 2021
$wsCustos.Range("D21").Value = 20

# column widths / bestFit for the new data
$wsCustos.Columns.Item(1).ColumnWidth = 15.140625
$wsCustos.Columns.Item(4).ColumnWidth = 15.140625

# ---------------------------------------------------------------------------
# 4) Sheet-view / selection bookkeeping to match the authored end state:
#    Cenarios selection moves, Custos gets zoomed + selected range, and
#    Parametros ends up as the active sheet/tab.
# ---------------------------------------------------------------------------
$wsCenarios = $wb.Worksheets.Item("Cenarios")
$wsCenarios.Activate()
$wsCenarios.Range("A13").Select()

$wsCustos.Activate()
$wsCustos.Application.ActiveWindow.Zoom = 115
$wsCustos.Range("D3:D6").Select()

$wsParam.Activate()
$wsParam.Range("F24").Select()
